$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, J, K, L, M, P

$ws.Range("D2").Value2 = 44414
$ws.Range("J2").Value2 = 500
$ws.Range("K2").Value2 = 31000
$ws.Range("L2").Value2 = 32000
$ws.Range("M2").Value2 = 31500
$ws.Range("P2").Value2 = 1260

$ws.Range("D3").Value2 = 44827
$ws.Range("J3").Value2 = 300
$ws.Range("K3").Value2 = 30000
$ws.Range("L3").Value2 = 31000
$ws.Range("M3").Value2 = 30500
$ws.Range("P3").Value2 = 1220
